$d = $word.ActiveDocument

# Locate the paragraph that contains "LOB1039: Física Experimental III (Requisito fraco)"
# and then remove the three paragraphs that immediately followed it:
#   - a blank paragraph
#   - "Ver no Jupiter Salvar em pdf Salvar em docx"
#   - "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"
$paras = $d.Paragraphs

$anchorIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Text -like "*LOB1039*Requisito fraco*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -ge 1) {
    # Delete the three following paragraphs, starting from the last one so
    # earlier indices stay valid while we work.
    $paras.Item($anchorIndex + 3).Range.Delete()
    $paras.Item($anchorIndex + 2).Range.Delete()
    $paras.Item($anchorIndex + 1).Range.Delete()
}
